$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like
# "0.0690" or "0.000220" keep their exact textual representation
# instead of being coerced into numbers by Excel.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11:D12").NumberFormat = "@"
$ws.Range("D14:D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34:D40").NumberFormat = "@"
$ws.Range("D42:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.608.87"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "2.267.40"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "230.14"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "63.40"
$ws.Range("E7").Value = "  +4.54%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "0.423"
$ws.Range("E9").Value = "  +3.89%  "
$ws.Range("E10").Value = "  +8.46%  "
$ws.Range("D11").Value = "57.27"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "25.83"
$ws.Range("E12").Value = "  +13.58%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "2.604.27"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "15.63"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "5.87"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("D17").Value = "0.812"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "2.266.00"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Value = "43.562.72"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  +4.10%  "
$ws.Range("D21").Value = "72.96"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").Value = "248.51"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +4.21%  "
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").Value = "171.15"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").Value = "20.48"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("E32").Value = "  +10.90%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "0.0690"
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("D35").Value = "5.10"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Value = "4.68"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +4.48%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "6.72"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "2.32"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").Value = "0.0246"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "10.65"
$ws.Range("E42").Value = "  +21.56%  "
$ws.Range("B43").Value = "TerraClassic"
$ws.Range("C43").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D43").Value = "0.000220"
$ws.Range("E43").Value = "  -10.47%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "8.34"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "4.46"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0963"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").Value = "97.06"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "1.473.01"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "16.79"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  +1.31%  "
